$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-13 Saturday" "2024-04-14 Sunday"

Replace-Text "818÷7=" "703÷8="
Replace-Text "440÷5=" "833÷8="
Replace-Text "737÷3=" "791÷6="
Replace-Text "545÷5=" "505÷6="
Replace-Text "358÷6=" "131÷9="

Replace-Text "427÷3=" "581÷2="
Replace-Text "104÷7=" "448÷2="
Replace-Text "402÷5=" "915÷6="
Replace-Text "557÷8=" "735÷4="
Replace-Text "604÷8=" "702÷6="

Replace-Text "107÷2=" "330÷4="
Replace-Text "482÷5=" "870÷6="
Replace-Text "540÷3=" "603÷3="
Replace-Text "382÷6=" "908÷3="
Replace-Text "628÷5=" "528÷2="

Replace-Text "846÷3=" "811÷5="
Replace-Text "978÷2=" "473÷3="
Replace-Text "183÷8=" "692÷8="
Replace-Text "198÷2=" "682÷8="
Replace-Text "364÷7=" "485÷4="

Replace-Text "824÷7=" "685÷7="
Replace-Text "938÷5=" "233÷9="
Replace-Text "502÷4=" "952÷9="
Replace-Text "609÷9=" "383÷6="
Replace-Text "917÷7=" "474÷8="

Write-Output "Replacements complete"
